# Fill in the Wednesday / Thursday / Friday schedule blocks with the same
# class pattern already present for Tuesday (rows 9-15), spreading it across
# rows 16-22 (Wednesday), 23-29 (Thursday) and 30-36 (Friday).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($base in 16, 23, 30) {
    $r0 = $base       # 08:00 row
    $r1 = $base + 1   # 10:00 row
    $r2 = $base + 2   # 12:00 row
    $r3 = $base + 3   # 14:00 row
    $r4 = $base + 4   # 16:00 row
    $r5 = $base + 5   # 18:00 row

    $ws.Range("C$r0").Value = "PL-204"
    $ws.Range("E$r0").Value = "AM-BT"
    $ws.Range("F$r0").Value = "Fizica-192"

    $ws.Range("C$r1").Value = "PC-205"

    $ws.Range("D$r2").Value = "ENG-BT"
    $ws.Range("E$r2").Value = "PL-204"

    $ws.Range("E$r3").Value = "PC-205"
    $ws.Range("F$r3").Value = "ENG-BT"

    $ws.Range("D$r4").Value = "PL-204"

    $ws.Range("D$r5").Value = "PC-205"
}

# Restore the active selection recorded in the saved workbook.
$ws.Range("I26").Select()
